# E11_David_Mugwaneza.docx — "Updated exercise11 solutions doc"
#
# The underlying commit is dominated by Word's background proofing
# engine silently wrapping already-correct text in <w:proofErr> marks
# (gramStart/gramEnd/spellStart/spellEnd) when the file was re-opened
# and re-saved in Word — the visible text itself is byte-for-byte
# identical before/after every one of those hunks. Those marks have no
# COM-automatable surface (no VBA/COM call ever inserted them even in
# real Word; they are an editor-internal side effect of interactive
# proofing), so this script focuses on the one concrete, user-visible,
# COM-reachable edit in the diff: the Alt Text ("descr") that was added
# to the last picture (id 1379440164, "Picture 2") via the
# Format Picture > Alt Text pane, which Word mirrors onto both the
# <wp:docPr> and <pic:cNvPr> nodes of that drawing.

$d = $word.ActiveDocument

$pic = $d.InlineShapes.Item(14)

$altText = "A black square with white dots&#10;&#10;Description automatically generated"
$pic.AlternativeText = $altText

Write-Output ("AlternativeText set to: " + $pic.AlternativeText)
